# Change the table style used by the three tables in this deck from the
# custom "Table Grid" style to the built-in "Medium Style 2 - Accent 1"
# style (Table Design gallery swap).

$p = $ppt.ActivePresentation

$oldStyleId = "{4A7E6A58-6B7A-4B7E-9E24-32520A5AABDC}"
$newStyleId = "{E68C4403-F217-4881-B336-E87157FC7A2C}"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)

    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)

        if ($shape.HasTable) {
            $table = $shape.Table

            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
